# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Price cells that look like plain decimals (e.g. "1.002") are written with a
# leading apostrophe so Excel keeps them as literal text instead of coercing
# them to numbers (which would silently drop trailing zeros / change value),
# then the cell style is reset back to "Normal" so no visible quote prefix or
# stray number-format remains - this matches the original text-cell layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.206.91'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.870.49'
$ws.Range("E3").Value = '  +2.02%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").Value = '''311.67'
$ws.Range("E5").Value = '  +0.24%  '
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("D7").Value = '''0.5052'
$ws.Range("E7").Value = '  -1.33%  '
$ws.Range("D8").Value = '''0.3916'
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").Value = '''0.09631'
$ws.Range("E9").Value = '  -5.78%  '
$ws.Range("D10").Value = '''1.138'
$ws.Range("E10").Value = '  +2.43%  '
$ws.Range("D11").Value = '''40.83'
$ws.Range("E11").Value = '  -0.25%  '
$ws.Range("D12").Value = '''6.491'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '''20.92'
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("D14").Value = '1.882.46'
$ws.Range("E14").Value = '  +4.03%  '
$ws.Range("D15").Value = '''7.434'
$ws.Range("E15").Value = '  +0.29%  '
$ws.Range("E16").Value = '  +0.41%  '
$ws.Range("D17").Value = '''0.00001125'
$ws.Range("E17").Value = '  -2.44%  '
$ws.Range("D18").Value = '''92.88'
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("D19").Value = '''0.06625'
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("D20").Value = '''17.55'
$ws.Range("E20").Value = '  +1.02%  '
$ws.Range("E21").Value = '  +0.25%  '
$ws.Range("D22").Value = '''6.146'
$ws.Range("E22").Value = '  +1.39%  '
$ws.Range("D23").Value = '28.260.54'
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("E24").Value = '  +1.56%  '
$ws.Range("E25").Value = '  +1.45%  '
$ws.Range("D26").Value = '''2.535'
$ws.Range("E26").Value = '  +3.57%  '
$ws.Range("D27").Value = '2.080.30'
$ws.Range("E27").Value = '  +2.42%  '
$ws.Range("D28").Value = '''21.19'
$ws.Range("E28").Value = '  +2.59%  '
$ws.Range("D29").Value = '''157.44'
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '''127.23'
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").Value = '''0.1059'
$ws.Range("E31").Value = '  -3.47%  '
$ws.Range("D32").Value = '''1.066'
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("D33").Value = '''5.623'
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("D34").Value = '''3.625'
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").Value = '''9.642'
$ws.Range("E35").Value = '  +6.20%  '
$ws.Range("E36").Value = '  -2.40%  '
$ws.Range("D37").Value = '''0.02386'
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '''11.46'
$ws.Range("E39").Value = '  -1.13%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.6345'
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("D41").Value = '''4.981'
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("D42").Value = '''1.177'
$ws.Range("E42").Value = '  +2.04%  '
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("D44").Value = '''13.51'
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("E45").Value = '  +0.43%  '
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = '''1.267'
$ws.Range("E47").Value = '  -1.86%  '
$ws.Range("D48").Value = '''124.62'
$ws.Range("E48").Value = '  -0.43%  '
$ws.Range("D49").Value = '''1.990'
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  +0.68%  '

$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D17", "D18", "D19", "D20", "D22", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39", "D40", "D41", "D42", "D44", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
